# Relative Wealth index prep and added across the dashboard
#
# Adds two new indicator rows (Mean / Majority Relative Wealth Index) to the
# legend_data sheet, right after the existing "Share of population Living in
# Area with Low Road Density" row, and moves the sheet's scroll/selection to
# the newly added data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 61: Mean Relative Wealth Index -----------------------------------
$ws.Range("A61").Value = "RWI_rwi_mean"
$ws.Range("B61").Value = "positive"
$ws.Range("D61").Value = "Pakistan Poverty Team - World Bank"
$ws.Range("E61").Value = "Mean Relative Wealth Index"
$ws.Range("F61").Value = "Mean Relative Wealth Index"

# --- Row 62: Majority Relative Wealth Index --------------------------------
$ws.Range("A62").Value = "RWI_rwi_majority"
$ws.Range("B62").Value = "positive"
$ws.Range("D62").Value = "Pakistan Poverty Team - World Bank"
$ws.Range("E62").Value = "Majority Relative Wealth Index"
$ws.Range("F62").Value = "Majority Relative Wealth Index"

# --- Move the viewport / selection to the newly added rows -----------------
$win = $excel.ActiveWindow
$win.ScrollRow = 46
$win.ScrollColumn = 5
$ws.Range("G61").Select() | Out-Null
